# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Siren_Profits leve-profit workbook (8 job sheets: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each row below corresponds to one <row> hunk in the source diff; cell values are
# taken verbatim from the diff's "+" side. A handful of cells are removed entirely
# (ClearContents) where the diff shows a <c> element disappearing rather than just its <v>.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 871.9091
$ws.Range("I2").Value = 909.1
$ws.Range("K2").Value = 909.1
$ws.Range("M2").Value = -796.1

# Row 34
$ws.Range("H34").Value = 1055.7273
$ws.Range("I34").Value = 845.8889
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 845.8889
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -642.8889
$ws.Range("N34").Value = -2406

# Row 36
$ws.Range("H36").Value = 1055.7273
$ws.Range("I36").Value = 845.8889
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 845.8889
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = -130.8889
$ws.Range("N36").Value = -3430

# Row 62
$ws.Range("H62").Value = 71458930
$ws.Range("I62").Value = 333343330
$ws.Range("K62").Value = 333343330
$ws.Range("M62").Value = -333342706

# Row 65
$ws.Range("H65").Value = 71458930
$ws.Range("I65").Value = 333343330
$ws.Range("K65").Value = 1666716650
$ws.Range("M65").Value = -1666713530

# Row 107
$ws.Range("H107").Value = 18750
$ws.Range("I107").Value = 21666.666
$ws.Range("K107").Value = 21666.666
$ws.Range("M107").Value = -19746.666

# Row 135
$ws.Range("H135").Value = 7665.6665
$ws.Range("I135").Value = 8798.799999999999
$ws.Range("K135").Value = 79189.2
$ws.Range("M135").Value = -76654.2

# Row 137
$ws.Range("H137").Value = 7734.4243
$ws.Range("I137").Value = 9526.32
$ws.Range("J137").Value = 2134.75
$ws.Range("K137").Value = 28578.96
$ws.Range("L137").Value = 6404.25
$ws.Range("M137").Value = -26028.96
$ws.Range("N137").Value = -11504.25

# Row 138
$ws.Range("H138").Value = 3922.4
$ws.Range("I138").Value = 629.8889
$ws.Range("J138").Value = 4645.1465
$ws.Range("K138").Value = 1889.6667
$ws.Range("L138").Value = 13935.4395
$ws.Range("M138").Value = 3250.3333
$ws.Range("N138").Value = -24215.4395

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1831.5186
$ws.Range("I32").Value = 1836.5577
$ws.Range("K32").Value = 1836.5577
$ws.Range("M32").Value = -1549.5577

# Row 45
$ws.Range("H45").Value = 6478.615
$ws.Range("I45").Value = 5284.125
$ws.Range("K45").Value = 5284.125
$ws.Range("M45").Value = -4907.125

# Row 97
$ws.Range("H97").Value = 19734.834
$ws.Range("I97").Value = 8411.857
$ws.Range("J97").Value = 59365.25
$ws.Range("K97").Value = 8411.857
$ws.Range("L97").Value = 59365.25
$ws.Range("M97").Value = -7915.857
$ws.Range("N97").Value = -60357.25

# Row 133
$ws.Range("H133").Value = 78000
$ws.Range("J133").Value = 78000
$ws.Range("L133").Value = 78000
$ws.Range("N133").Value = -83060

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 892
$ws.Range("I80").Value = 522.5
$ws.Range("J80").Value = 948.8461
$ws.Range("K80").Value = 522.5
$ws.Range("L80").Value = 948.8461
$ws.Range("M80").Value = 475.5
$ws.Range("N80").Value = -2944.8461

# Row 83
$ws.Range("H83").Value = 892
$ws.Range("I83").Value = 522.5
$ws.Range("J83").Value = 948.8461
$ws.Range("K83").Value = 2612.5
$ws.Range("L83").Value = 4744.2305
$ws.Range("M83").Value = 2379.5
$ws.Range("N83").Value = -14728.2305

# Row 94
$ws.Range("H94").Value = 2957.9333
$ws.Range("I94").Value = 3228.2222
$ws.Range("J94").Value = 2552.5
$ws.Range("K94").Value = 3228.2222
$ws.Range("L94").Value = 2552.5
$ws.Range("M94").Value = -2777.2222
$ws.Range("N94").Value = -3454.5

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 18363.334
$ws.Range("I51").Value = 10045
$ws.Range("J51").Value = 35000
$ws.Range("K51").Value = 10045
$ws.Range("L51").Value = 35000
$ws.Range("M51").Value = -9309
$ws.Range("N51").Value = -36472

# Row 58
$ws.Range("H58").Value = 1714.3043
$ws.Range("I58").Value = 1229.5
$ws.Range("J58").Value = 3459.6
$ws.Range("K58").Value = 1229.5
$ws.Range("L58").Value = 3459.6
$ws.Range("M58").Value = -1026.5
$ws.Range("N58").Value = -3865.6

# Row 61
$ws.Range("H61").Value = 18363.334
$ws.Range("I61").Value = 10045
$ws.Range("J61").Value = 35000
$ws.Range("K61").Value = 10045
$ws.Range("L61").Value = 35000
$ws.Range("M61").Value = -9697
$ws.Range("N61").Value = -35696

# Row 107
$ws.Range("H107").Value = 13724.294
$ws.Range("I107").Value = 20389.727
$ws.Range("J107").Value = 1504.3334
$ws.Range("K107").Value = 20389.727
$ws.Range("L107").Value = 1504.3334
$ws.Range("M107").Value = -18469.727
$ws.Range("N107").Value = -5344.3334

# Row 136
$ws.Range("H136").Value = 1714.3043
$ws.Range("I136").Value = 1229.5
$ws.Range("J136").Value = 3459.6
$ws.Range("K136").Value = 3688.5
$ws.Range("L136").Value = 10378.8
$ws.Range("M136").Value = -1138.5
$ws.Range("N136").Value = -15478.8

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 47795.273
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 27000000
$ws.Range("I18").Value = 27000000
$ws.Range("K18").Value = 27000000
$ws.Range("M18").Value = -26999707

# Row 97
$ws.Range("H97").Value = 14045.235
$ws.Range("I97").Value = 19823.727
$ws.Range("J97").Value = 3451.3333
$ws.Range("K97").Value = 19823.727
$ws.Range("L97").Value = 3451.3333
$ws.Range("M97").Value = -19327.727
$ws.Range("N97").Value = -4443.3333

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 13307.846
$ws.Range("I61").Value = 15833.833
$ws.Range("K61").Value = 15833.833
$ws.Range("M61").Value = -15631.833

# Row 93
$ws.Range("H93").Value = 16829.666
$ws.Range("I93").Value = 16829.666
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 16829.666
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -15581.666
$ws.Range("N93").ClearContents()

# Row 100
$ws.Range("H100").Value = 2458.1667
$ws.Range("I100").Value = 2083
$ws.Range("J100").Value = 2833.3333
$ws.Range("K100").Value = 2083
$ws.Range("L100").Value = 2833.3333
$ws.Range("M100").Value = -1542
$ws.Range("N100").Value = -3915.3333

# Row 113
$ws.Range("H113").Value = 13307.846
$ws.Range("I113").Value = 15833.833
$ws.Range("K113").Value = 15833.833
$ws.Range("M113").Value = -13663.833

$ws = $wb.Worksheets.Item("WVR")
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 132
$ws.Range("H132").Value = 11680.071
$ws.Range("I132").Value = 13050.419
$ws.Range("K132").Value = 39151.257
$ws.Range("M132").Value = -36621.257

Write-Host "Applied market-data refresh to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"